# Remove the rows for the accounts that were dropped from the "Export" sheet.
# (Unified diff shows 35 data rows removed from xl/worksheets/sheet.xml while
# everything else - headers, remaining rows, footer note - is unchanged.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$accountsToRemove = @(
    "004328934",
    "004486497",
    "005366255",
    "004352384",
    "004500804",
    "005079458",
    "005152037",
    "004321092",
    "004556853",
    "003435941",
    "004946997",
    "005055239",
    "004971783",
    "005256036",
    "004854514",
    "004983378",
    "004369172",
    "004376145",
    "001759765",
    "004565108",
    "004267044",
    "004693308",
    "004973881",
    "005654767",
    "005701765",
    "004927044",
    "005312963",
    "004547722",
    "004453045",
    "004436055",
    "004212476",
    "004983395",
    "004212581",
    "005186167",
    "004862677"
)

$lastRow = $ws.UsedRange.Rows.Count
$deletedCount = 0

# Walk bottom-up so deleting a row doesn't shift the index of rows not yet visited.
for ($r = $lastRow; $r -ge 2; $r--) {
    $acct = $ws.Cells.Item($r, 1).Value()
    if ($accountsToRemove -contains $acct) {
        $ws.Cells.Item($r, 1).EntireRow.Delete()
        $deletedCount = $deletedCount + 1
    }
}

Write-Output "Deleted $deletedCount rows"
